$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (B13/C13 = "984972 - Hugo Ricardo Zschommler Sandim") is removed
# entirely; everything below shifts up by one row.
$ws.Rows("13:13").Delete()

# After the shift, apply the content updates for the rows that actually changed
# text (row numbers below are POST-shift, i.e. final row numbers).

# Row 13 ("Programa resumido:") - summary replaced with "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:") - long syllabus text replaced with a date string
$ws.Range("B15").Value = "01/01/2016"
$ws.Range("C15").Value = "01/01/2016"

# Row 18 ("Método:") - now holds the professor's name
$ws.Range("B18").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C18").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# Row 19 ("Critério:") - now holds the evaluation method text
$ws.Range("B19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."
$ws.Range("C19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."

# Row 20 ("Norma de recuperação:") - now holds the final-grade formula
$ws.Range("B20").Value = "Nota Final NF = [P1 + P2]/2"
$ws.Range("C20").Value = "Nota Final NF = [P1 + P2]/2"

# Row 21 ("Bibliografia:") - now holds the recovery-exam rule text
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"
